# Apply the model updates to the "model" sheet of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# --- Core assumption changes ---

# Year-5 (V-column) revenue assumption bumped from 37,000 to 40,000.
# (W39:Z39 are formulas referencing V39, so they recalc automatically.)
$ws.Range("V39").Value = 40000

# Margin ramp (row 59): V59 nudged up, and W59/X59 become a formula-driven
# step-up (prior year * 1.02) instead of hardcoded values; Y59/Z59 settle
# at a lower hardcoded 0.13 instead of 0.14.
$ws.Range("V59").Value = 0.125
$ws.Range("W59").Formula = "=V59*1.02"
$ws.Range("X59").Formula = "=W59*1.02"
$ws.Range("Y59").Value = 0.13
$ws.Range("Z59").Value = 0.13

# --- View-state bookkeeping (best-effort; matches where the author's
# selection ended up after the edits) ---
$ws.Select()
$ws.Range("AC56").Select()
